$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of J2:J11 (k column)
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# Row 14-17: summary statistics with labels
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Build a template style (bold, size 12, vertical-center) on a scratch cell,
# then paste just the formatting onto B14:B17 in a single operation so the
# style table only gains one new entry (matches native Excel behaviour).
$scratch = $ws.Range("AB1")
$scratch.Value = "x"
$scratch.Font.Bold = $true
$scratch.Font.Size = 12
$scratch.VerticalAlignment = -4108  # xlVAlignCenter

$scratch.Copy() | Out-Null
$ws.Range("B14:B17").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$scratch.Clear() | Out-Null

# Row heights for the new summary rows
$ws.Range("A14:A17").RowHeight = 15.6

# Page setup (paper size / orientation) as captured in the saved file
$ws.PageSetup.PaperSize = 9      # xlPaperA4
$ws.PageSetup.Orientation = 1    # xlPortrait

# Selection as in diff
$ws.Range("A14:B17").Select()
